$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking values (e.g. "1.001") are
# preserved as literal text instead of being parsed into floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "25.821.31"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "1.738.97"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "231.57"
$ws.Range("E5").Value = "  -1.76%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").Value = "0.5165"
$ws.Range("E7").Value = "  +1.66%  "
$ws.Range("D8").Value = "0.2805"
$ws.Range("E8").Value = "  +5.04%  "
$ws.Range("E9").Value = "  -2.91%  "
$ws.Range("D10").Value = "0.06115"
$ws.Range("D11").Value = "1.750.84"
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("D12").Value = "0.07035"
$ws.Range("E12").Value = "  +1.28%  "
$ws.Range("D13").Value = "15.29"
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("D14").Value = "0.6474"
$ws.Range("E14").Value = "  +4.26%  "
$ws.Range("D15").Value = "4.526"
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("D16").Value = "76.97"
$ws.Range("E16").Value = "  -0.92%  "
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").Value = "1.0000"
$ws.Range("E18").Value = "  -0.08%  "
$ws.Range("D19").Value = "25.819.91"
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("D20").Value = "11.48"
$ws.Range("E20").Value = "  -0.80%  "
$ws.Range("D21").Value = "0.000006600"
$ws.Range("E21").Value = "  -0.47%  "
$ws.Range("D22").Value = "1.975.58"
$ws.Range("E22").Value = "  +0.40%  "
$ws.Range("D23").Value = "4.138"
$ws.Range("E23").Value = "  +2.34%  "
$ws.Range("D24").Value = "8.661"
$ws.Range("E24").Value = "  +5.02%  "
$ws.Range("D25").Value = "5.141"
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("D26").Value = "139.31"
$ws.Range("E26").Value = "  +2.06%  "
$ws.Range("D27").Value = "1.510"
$ws.Range("E27").Value = "  +3.32%  "
$ws.Range("D28").Value = "15.07"
$ws.Range("E28").Value = "  +0.59%  "
$ws.Range("D29").Value = "1.814"
$ws.Range("E29").Value = "  +2.83%  "
$ws.Range("D30").Value = "102.17"
$ws.Range("E30").Value = "  -0.24%  "
$ws.Range("D31").Value = "0.08316"
$ws.Range("E31").Value = "  +2.30%  "
$ws.Range("D32").Value = "3.680"
$ws.Range("E32").Value = "  -0.16%  "
$ws.Range("D33").Value = "3.428"
$ws.Range("E33").Value = "  +1.31%  "
$ws.Range("D34").Value = "0.04494"
$ws.Range("E34").Value = "  +2.22%  "
$ws.Range("D35").Value = "2.609"
$ws.Range("E35").Value = "  -1.54%  "
$ws.Range("D36").Value = "0.9855"
$ws.Range("E36").Value = "  -0.73%  "
$ws.Range("D37").Value = "0.6148"
$ws.Range("E37").Value = "  +2.18%  "
$ws.Range("D38").Value = "2.658"
$ws.Range("E38").Value = "  +2.46%  "
$ws.Range("E39").Value = "  +1.94%  "
$ws.Range("D40").Value = "1.938"
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("D41").Value = "0.9995"
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("D42").Value = "100.75"
$ws.Range("E42").Value = "  -0.78%  "
$ws.Range("D43").Value = "0.3846"
$ws.Range("E43").Value = "  +0.79%  "
$ws.Range("D44").Value = "0.7276"
$ws.Range("E44").Value = "  -2.40%  "
$ws.Range("D45").Value = "4.974"
$ws.Range("E45").Value = "  +1.57%  "
$ws.Range("D46").Value = "0.05410"
$ws.Range("E46").Value = "  -1.60%  "
$ws.Range("D47").Value = "6.280"
$ws.Range("E47").Value = "  +6.27%  "
$ws.Range("D48").Value = "0.1121"
$ws.Range("E48").Value = "  +2.46%  "
$ws.Range("D49").Value = "53.16"
$ws.Range("E49").Value = "  +1.41%  "
$ws.Range("D50").Value = "7.703"
$ws.Range("E50").Value = "  +4.18%  "
$ws.Range("D51").Value = "29.95"
$ws.Range("E51").Value = "  -0.19%  "

# Restore default (Normal) style on column D so no stray number format is left
# applied to the cells (matches original workbook formatting).
$ws.Range("D2:D51").Style = "Normal"

